$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.183.78"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "3.801.76"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'601.81"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "'165.26"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("E11").Value = "  +2.95%  "
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("D13").Value = "'35.89"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").Value = "4.441.08"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").Value = "3.811.37"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").Value = "68.193.55"
$ws.Range("D17").Value = "'18.49"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'461.88"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'9.73"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'0.0000149"
$ws.Range("E23").Value = "  -3.30%  "
$ws.Range("D24").Value = "'83.16"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "'10.01"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Value = "3.950.81"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -5.12%  "
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").Value = "'7.35"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").Value = "'29.43"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D36").Value = "'0.1000"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  -3.24%  "
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").Value = "'5.85"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").Value = "'0.989"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "'0.300"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").Value = "'47.55"
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").Value = "'43.24"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("D46").Value = "'151.51"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").Value = "'8.38"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").Value = "'1.87"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("D49").Value = "'397.51"
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("E50").Value = "  +4.45%  "
$ws.Range("D51").Value = "'26.76"
$ws.Range("E51").Value = "  +0.14%  "
